# Fill in the "NIP Wali Kelas" (F) and "Nama Wali Kelas" (G) columns for
# every class row. The F-column NIPs are long, zero-padded digit strings
# that must stay as text (leading zeros preserved), so the range is
# pre-formatted as Text before the values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F10").NumberFormat = "@"

$ws.Range("F2").Value = "0000000000000066"
$ws.Range("G2").Value = "M. Fais Jainuddin, S.Pd"

$ws.Range("F3").Value = "5736762663300210"
$ws.Range("G3").Value = "Nunung Indrawati, S.Pd."

$ws.Range("F4").Value = "0000000000000044"
$ws.Range("G4").Value = "Zulfi Amaliyah, S.Kom."

$ws.Range("F5").Value = "0000000000000006"
$ws.Range("G5").Value = "Mulyono, S.Th."

$ws.Range("F6").Value = "5040758659300040"
$ws.Range("G6").Value = "Nurmala Evayanti S.Pd."

$ws.Range("F7").Value = "00000000000000022222"
$ws.Range("G7").Value = "Fera Mega Haristina, S.Tr.Kom."

$ws.Range("F8").Value = "00000000000000004"
$ws.Range("G8").Value = "Imtiana, S.Pd."

$ws.Range("F9").Value = "0000000023232323"
$ws.Range("G9").Value = "Frances Laurence Setyo Budi, S.Pd."

$ws.Range("F10").Value = "0000000000000010044"
$ws.Range("G10").Value = "Udayani, S.Pd."
